$d = $word.ActiveDocument

# wdColor value for RGB C9211E (Word stores OLE colors as 0x00BBGGRR).
$commentColor = 1974729

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.Contains("I was searching a bit") -and $t.Contains("add location")) {
        $pr = $p.Range
        # End of the paragraph's text content, just before the paragraph mark.
        $insertAt = $pr.End - 1

        $rAll = $d.Range($insertAt, $insertAt)
        $rAll.InsertAfter("– Done")

        # Split the inserted text into its two runs ("– " and "Done"), each
        # carrying the reviewer-comment red color, matching the existing
        # "– Done" annotations used elsewhere in this document.
        $dashLen = 2
        $r1 = $d.Range($rAll.Start, $rAll.Start + $dashLen)
        $r1.Font.Color = $commentColor

        $r2 = $d.Range($rAll.Start + $dashLen, $rAll.End)
        $r2.Font.Color = $commentColor

        break
    }
}
